$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Load" sheet: PV distributed-generation rows 32-45, column C
#    (rated/base line-to-neutral voltage) updated from 2.4 to 4.1569219
#    to reflect the GMDM ePHASORSIM profile update.
# ---------------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Load")
$wsLoad.Range("C32:C45").Value = 4.1569219

# ---------------------------------------------------------------------------
# 2) "Bus" sheet: fixed GLM tank phasing. For several buses the phase
#    labels (column A) and angles (column E) were rotated one row within
#    their A/B/C (or A/C) group to correct the phase assignment.
#    Each tuple is: (row, new bus name, new angle).
# ---------------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("Bus")

$busData = @(
    @(2, "1_A", 0),
    @(3, "1_B", -120),
    @(4, "1_C", 120),
    @(6, "100_A", 0),
    @(7, "100_B", -120),
    @(8, "100_C", 120),
    @(9, "101_A", 0),
    @(10, "101_B", -120),
    @(11, "101_C", 120),
    @(15, "105_A", 0),
    @(16, "105_B", -120),
    @(17, "105_C", 120),
    @(20, "108_A", 0),
    @(21, "108_B", -120),
    @(22, "108_C", 120),
    @(31, "13_A", 0),
    @(32, "13_B", -120),
    @(33, "13_C", 120),
    @(34, "135_A", 0),
    @(35, "135_B", -120),
    @(36, "135_C", 120),
    @(38, "149_A", 0),
    @(39, "149_B", -120),
    @(40, "149_C", 120),
    @(42, "150_A", 0),
    @(43, "150_B", -120),
    @(44, "150_C", 120),
    @(45, "150r_A", 0),
    @(46, "150r_B", -120),
    @(47, "150r_C", 120),
    @(48, "151_A", 0),
    @(49, "151_B", -120),
    @(50, "151_C", 120),
    @(51, "152_A", 0),
    @(52, "152_B", -120),
    @(53, "152_C", 120),
    @(55, "160_A", 0),
    @(56, "160_B", -120),
    @(57, "160_C", 120),
    @(58, "160r_A", 0),
    @(59, "160r_B", -120),
    @(60, "160r_C", 120),
    @(62, "18_A", 0),
    @(63, "18_B", -120),
    @(64, "18_C", 120),
    @(66, "197_A", 0),
    @(67, "197_B", -120),
    @(68, "197_C", 120),
    @(71, "21_A", 0),
    @(72, "21_B", -120),
    @(73, "21_C", 120),
    @(75, "23_A", 0),
    @(76, "23_B", -120),
    @(77, "23_C", 120),
    @(79, "25_A", 0),
    @(80, "25_B", -120),
    @(81, "25_C", 120),
    @(82, "250_A", 0),
    @(83, "250_B", -120),
    @(84, "250_C", 120),
    @(85, "25r_A", 0),
    @(86, "25r_C", 120),
    @(87, "26_A", 0),
    @(88, "26_C", 120),
    @(89, "27_A", 0),
    @(90, "27_C", 120),
    @(91, "28_A", 0),
    @(92, "28_B", -120),
    @(93, "28_C", 120),
    @(94, "29_A", 0),
    @(95, "29_B", -120),
    @(96, "29_C", 120),
    @(98, "30_A", 0),
    @(99, "30_B", -120),
    @(100, "30_C", 120),
    @(101, "300_A", 0),
    @(102, "300_B", -120),
    @(103, "300_C", 120),
    @(104, "300_open_A", 0),
    @(105, "300_open_B", -120),
    @(106, "300_open_C", 120),
    @(111, "35_A", 0),
    @(112, "35_B", -120),
    @(113, "35_C", 120),
    @(114, "36_A", 0),
    @(115, "36_B", -120),
    @(120, "40_A", 0),
    @(121, "40_B", -120),
    @(122, "40_C", 120),
    @(124, "42_A", 0),
    @(125, "42_B", -120),
    @(126, "42_C", 120),
    @(128, "44_A", 0),
    @(129, "44_B", -120),
    @(130, "44_C", 120),
    @(132, "450_A", 0),
    @(133, "450_B", -120),
    @(134, "450_C", 120),
    @(136, "47_A", 0),
    @(137, "47_B", -120),
    @(138, "47_C", 120),
    @(139, "48_A", 0),
    @(140, "48_B", -120),
    @(141, "48_C", 120),
    @(142, "49_A", 0),
    @(143, "49_B", -120),
    @(144, "49_C", 120),
    @(146, "50_A", 0),
    @(147, "50_B", -120),
    @(148, "50_C", 120),
    @(149, "51_A", 0),
    @(150, "51_B", -120),
    @(151, "51_C", 120),
    @(152, "52_A", 0),
    @(153, "52_B", -120),
    @(154, "52_C", 120),
    @(155, "53_A", 0),
    @(156, "53_B", -120),
    @(157, "53_C", 120),
    @(158, "54_A", 0),
    @(159, "54_B", -120),
    @(160, "54_C", 120),
    @(161, "55_A", 0),
    @(162, "55_B", -120),
    @(163, "55_C", 120),
    @(164, "56_A", 0),
    @(165, "56_B", -120),
    @(166, "56_C", 120),
    @(167, "57_A", 0),
    @(168, "57_B", -120),
    @(169, "57_C", 120),
    @(173, "60_A", 0),
    @(174, "60_B", -120),
    @(175, "60_C", 120),
    @(176, "61_A", 0),
    @(177, "61_B", -120),
    @(178, "61_C", 120),
    @(179, "610_A", 0),
    @(180, "610_B", -120),
    @(181, "610_C", 120),
    @(182, "61s_A", 0),
    @(183, "61s_B", -120),
    @(184, "61s_C", 120),
    @(185, "62_A", 0),
    @(186, "62_B", -120),
    @(187, "62_C", 120),
    @(188, "63_A", 0),
    @(189, "63_B", -120),
    @(190, "63_C", 120),
    @(191, "64_A", 0),
    @(192, "64_B", -120),
    @(193, "64_C", 120),
    @(194, "65_A", 0),
    @(195, "65_B", -120),
    @(196, "65_C", 120),
    @(197, "66_A", 0),
    @(198, "66_B", -120),
    @(199, "66_C", 120),
    @(200, "67_A", 0),
    @(201, "67_B", -120),
    @(202, "67_C", 120),
    @(205, "7_A", 0),
    @(206, "7_B", -120),
    @(207, "7_C", 120),
    @(210, "72_A", 0),
    @(211, "72_B", -120),
    @(212, "72_C", 120),
    @(216, "76_A", 0),
    @(217, "76_B", -120),
    @(218, "76_C", 120),
    @(219, "77_A", 0),
    @(220, "77_B", -120),
    @(221, "77_C", 120),
    @(222, "78_A", 0),
    @(223, "78_B", -120),
    @(224, "78_C", 120),
    @(225, "79_A", 0),
    @(226, "79_B", -120),
    @(227, "79_C", 120),
    @(228, "8_A", 0),
    @(229, "8_B", -120),
    @(230, "8_C", 120),
    @(231, "80_A", 0),
    @(232, "80_B", -120),
    @(233, "80_C", 120),
    @(234, "81_A", 0),
    @(235, "81_B", -120),
    @(236, "81_C", 120),
    @(237, "82_A", 0),
    @(238, "82_B", -120),
    @(239, "82_C", 120),
    @(240, "83_A", 0),
    @(241, "83_B", -120),
    @(242, "83_C", 120),
    @(245, "86_A", 0),
    @(246, "86_B", -120),
    @(247, "86_C", 120),
    @(248, "87_A", 0),
    @(249, "87_B", -120),
    @(250, "87_C", 120),
    @(252, "89_A", 0),
    @(253, "89_B", -120),
    @(254, "89_C", 120),
    @(257, "91_A", 0),
    @(258, "91_B", -120),
    @(259, "91_C", 120),
    @(261, "93_A", 0),
    @(262, "93_B", -120),
    @(263, "93_C", 120),
    @(266, "95_A", 0),
    @(267, "95_B", -120),
    @(268, "95_C", 120),
    @(270, "97_A", 0),
    @(271, "97_B", -120),
    @(272, "97_C", 120),
    @(273, "98_A", 0),
    @(274, "98_B", -120),
    @(275, "98_C", 120),
    @(276, "99_A", 0),
    @(277, "99_B", -120),
    @(278, "99_C", 120)
)

foreach ($item in $busData) {
    $r = $item[0]
    $name = $item[1]
    $ang = $item[2]
    $wsBus.Cells.Item($r, 1).Value = $name
    $wsBus.Cells.Item($r, 5).Value = $ang
}
